# Add new columns I (I0) and J (IF) to the worksheet, mirroring the
# existing header style (s="1", bold/bordered) for the header cells and
# plain numeric cells for the data rows (rows 2-60).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 / J1 - copy style from existing header cell H1 so the
# new headers get the same bold/border/centered formatting.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$iValues = @(8,8,7,8,6,8,7,8,7,7,7,7,7,9,6,8,8,8,6,8,7,7,7,8,9,8,8,5,8,8,7,7,6,8,7,6,7,7,7,7,7,8,7,7,8,7,7,6,9,7,5,10,6,7,7,7,5,6,9)
$jValues = @(8,8,7,8,6,8,7,8,8,8,7,7,8,9,6,8,8,8,7,8,8,7,8,8,9,8,8,6,8,8,7,7,7,9,7,7,8,8,7,7,7,9,8,7,8,8,7,6,9,7,6,10,6,7,8,8,5,6,9)

for ($r = 2; $r -le 60; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
